# Populate the "downloads" column (I) with a Python-dict-style string that
# wraps the existing "Data sheets" URL (column H) for each product row that
# currently holds the placeholder "(nan, nan)".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 193; $row -le 230; $row++) {
    $dataSheetUrl = $ws.Cells.Item($row, 8).Text   # column H = "Data sheets"
    $ws.Cells.Item($row, 9).Value = "{'Data sheet': ['" + $dataSheetUrl + "']}"
}
